$wb = $excel.ActiveWorkbook

# Sheet "展览" (first sheet): update want-to-go counts for rows 5 and 6
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 2927
$ws1.Range("F6").Value = 293

# Sheet "全部类型" (last sheet): same update, mirrored data
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 2927
$ws4.Range("F6").Value = 293
